$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 3 futures data (appended below the existing Week 1 / Week 2 blocks,
# which occupy rows 2-51 and 52-101 respectively).
$rows = @(
  @(280, "Georgia", 3),
  @(400, "Ohio State", 3),
  @(500, "Texas", 3),
  @(1800, "Penn State", 3),
  @(1100, "Oregon", 3),
  @(1200, "Alabama", 3),
  @(1400, "Ole Miss", 3),
  @(2500, "Tennessee", 3),
  @(2000, "Miami", 3),
  @(4000, "Missouri", 3),
  @(4000, "USC", 3),
  @(5000, "Utah", 3),
  @(4000, "LSU", 3),
  @(6600, "Oklahoma", 3),
  @(6000, "Clemson", 3),
  @(12500, "Michigan", 3),
  @(8000, "Notre Dame", 3),
  @(8000, "Kansas State", 3),
  @(8000, "Texas A&M", 3),
  @(15000, "Nebraska", 3),
  @(12500, "Louisville", 3),
  @(20000, "Oklahoma State", 3),
  @(20000, "Iowa State", 3),
  @(25000, "Arizona", 3),
  @(40000, "Iowa", 3),
  @(20000, "Boise State", 3),
  @(25000, "Wisconsin", 3),
  @(30000, "UCF", 3),
  @(30000, "Kansas", 3),
  @(50000, "Texas Tech", 3),
  @(30000, "Syracuse", 3),
  @(50000, "Colorado", 3),
  @(30000, "Washington", 3),
  @(50000, "Oregon State", 3),
  @(25000, "NC State", 3),
  @(25000, "Auburn", 3),
  @(40000, "Virginia Tech", 3),
  @(40000, "North Carolina", 3),
  @(40000, "South Carolina", 3),
  @(100000, "Liberty", 3),
  @(60000, "Tulane", 3),
  @(40000, "West Virginia", 3),
  @(40000, "Florida State", 3),
  @(50000, "SMU", 3),
  @(100000, "Michigan State", 3),
  @(50000, "Georgia Tech", 3),
  @(50000, "Kentucky", 3),
  @(100000, "Baylor", 3),
  @(50000, "Maryland", 3),
  @(50000, "Memphis", 3)
)

$r = 102
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# Move the view: selection on D16, no frozen/scrolled top-left cell override.
$ws.Range("D16").Select()
